$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B44: was stored as text "2", should become a real number 2
$ws.Range("B44").Value = 2

# Add new row 45 with the new annotation entry
$ws.Range("A45").Value = "Sunsi Wu"

# B45 mirrors the historical quirk of this sheet (politeness_score stored as text)
$ws.Range("B45").Value = "'2"

$ws.Range("C45").Value = "very limited"
$ws.Range("D45").Value = "CRT"
$ws.Range("E45").Value = "MET"
$ws.Range("F45").Value = "17635cfa-5d3f-4715-99a5-e710de1fbea7"
$ws.Range("G45").Value = "S1XXq6lRW_annotated.xlsx"
$ws.Range("H45").Value = "Technical contribution of the paper is very limited."
